$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(85, 1).NumberFormat = $ws.Cells.Item(84, 1).NumberFormat

$ws.Cells.Item(85, 1).Value = 46034
$ws.Cells.Item(85, 2).Value = 192
$ws.Cells.Item(85, 3).Value = 205
$ws.Cells.Item(85, 4).Value = 191
